$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.3390307724475861, 0.9012500047683716),
    @(0.3807833194732666, 0.8287500143051147),
    @(0.3720756471157074, 0.8849999904632568),
    @(0.3701992034912109, 0.8887500166893005),
    @(0.3658941686153412, 0.8912500143051147)
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
